# Update p-values in the MDMA column of the "Adverse events x molecules"
# table (results/paper_tables/session/02_dr_by_ae_molecule_matrix.docx).
#
# Table layout: column 1 = Adverse event, column 2 = LSD, column 3 = PSILOCYBIN,
# column 4 = MDMA. Row 1 is the header row.
#
# Each entry below is: (row index, adverse event label [for sanity logging],
#                        old text, new text).
#
# Note: several p-value strings are duplicated elsewhere in the table (e.g.
# "* (p=0.049)" also appears for "illusion"/LSD, " (p=0.072)" also appears for
# "autonomic"/MDMA, " (p=0.372)" also appears for "myalgia"/PSILOCYBIN), so
# each replacement is scoped to the specific table cell rather than done as a
# document-wide Find/Replace. We locate the match inside the cell range first
# (without replacing) and then assign the new text to that located sub-range,
# which also preserves the run's xml:space="preserve" formatting.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$changes = @(
    @(10, "anxiety",                 " (p=0.242)",  " (p=0.315)"),
    @(14, "decreased concentration",  " (p=0.521)",  " (p=0.373)"),
    @(19, "dizziness",                "* (p=0.049)", "* (p=0.017)"),
    @(22, "fatigue",                  " (p=0.092)",  " (p=0.103)"),
    @(29, "irritability",             " (p=0.527)",  " (p=0.532)"),
    @(30, "jaw tension",              " (p=0.177)",  " (p=0.120)"),
    @(31, "lack of appetite",         " (p=0.064)",  "* (p=0.049)"),
    @(33, "muscle tension",           " (p=0.287)",  " (p=0.241)"),
    @(39, "perspiration",             " (p=0.072)",  " (p=0.065)"),
    @(40, "restlessness",             " (p=0.595)",  " (p=0.558)"),
    @(41, "rumination",               " (p=0.499)",  " (p=0.633)"),
    @(42, "sleep disorder",           " (p=0.372)",  " (p=0.385)"),
    @(46, "weakness",                 " (p=0.677)",  " (p=0.837)")
)

foreach ($change in $changes) {
    $rowIndex = $change[0]
    $label = $change[1]
    $old = $change[2]
    $new = $change[3]

    $cell = $t.Cell($rowIndex, 4)

    # Defensive check: make sure we are editing the row we think we are.
    $labelCell = $t.Cell($rowIndex, 1)
    $labelText = $labelCell.Range.Text
    if (-not $labelText.StartsWith($label)) {
        Write-Host "WARNING: row $rowIndex label mismatch. Expected '$label', found '$labelText'"
    }

    $cellRange = $cell.Range
    $matchRange = $cellRange.Duplicate
    $found = $matchRange.Find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        Write-Host "WARNING: could not find '$old' in row $rowIndex ($label)"
    } else {
        $matchRange.Text = $new
        Write-Host "Updated row $rowIndex ($label): '$old' -> '$new'"
    }
}
